# "le monstre et la fille"
# Add 3 new rating columns (regard en arrierre, regard en avant, posthumain)
# and 5 new songs to the "Tableau1" table on Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. New table columns K, L, M
# ---------------------------------------------------------------------------
$null = $lo.ListColumns.Add()
$ws.Range("K1").Value = "regard en arrierre"

$null = $lo.ListColumns.Add()
$ws.Range("L1").Value = "regard en avant"

$null = $lo.ListColumns.Add()
$ws.Range("M1").Value = "posthumain"

# ---------------------------------------------------------------------------
# 2. Fill in the new rating values for existing songs
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

$ws.Range("M19").Value = 5

$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3
$ws.Range("M20").Value = 4

$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 2

$ws.Range("K38").Value = 4
$ws.Range("L38").Value = 0

$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 1

$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 3

# ---------------------------------------------------------------------------
# 3. Add the 4 new table rows (row 46 already exists as a blank table row)
# ---------------------------------------------------------------------------
$null = $lo.ListRows.Add()
$null = $lo.ListRows.Add()
$null = $lo.ListRows.Add()
$null = $lo.ListRows.Add()

# Row 46 - final conclusion
$ws.Range("A46").Value = "final conclusion"
$ws.Hyperlinks.Add($ws.Range("B46"), "https://www.youtube.com/watch?v=TYevEJ6y7Uo")
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = 3
$ws.Range("F46").Value = 2
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3
$ws.Range("L46").Value = 3

# Row 47 - Epilogue : beginning
$ws.Range("A47").Value = "Epilogue : beginning"
$ws.Hyperlinks.Add($ws.Range("B47"), "https://www.youtube.com/watch?v=08XIghnIjWs")
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 2
$ws.Range("E47").Value = 2
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 3
$ws.Range("I47").Value = 2
$ws.Range("J47").Value = 4
$ws.Range("K47").Value = 2
$ws.Range("L47").Value = 5

# Row 48 - Faded snapshots and forgotten dreams
$ws.Range("A48").Value = "Faded snapshots and forgotten dreams"
$ws.Hyperlinks.Add($ws.Range("B48"), "https://www.youtube.com/watch?v=Tzu8gBR0joY")
$ws.Range("C48").Value = 2
$ws.Range("D48").Value = 4
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 1
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3
$ws.Range("L48").Value = 0

# Row 49 - One dark mare
$ws.Range("A49").Value = "One dark mare"
$ws.Hyperlinks.Add($ws.Range("B49"), "https://www.youtube.com/watch?v=GfXAojf56wg")
$ws.Range("C49").Value = 5
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 5
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 5

# Row 50 - Sophia
$ws.Range("A50").Value = "Sophia"
$ws.Hyperlinks.Add($ws.Range("B50"), "https://www.youtube.com/watch?v=cOSAl26hyBQ")
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 1
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 3
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 1
$ws.Range("L50").Value = 3
$ws.Range("M50").Value = 0

# ---------------------------------------------------------------------------
# 4. Column widths (best achievable approximations through ColumnWidth)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 33.5
$ws.Columns.Item(2).ColumnWidth = 6.16666666666667
$ws.Columns.Item(11).ColumnWidth = 18
$ws.Columns.Item(12).ColumnWidth = 15.8333333333333
$ws.Columns.Item(13).ColumnWidth = 10.6666666666667

# ---------------------------------------------------------------------------
# 5. Restore view state (selection / scroll position)
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 31
$ws.Range("M43").Select()
